# Updates the cryptos list worksheet (prices + 1h volume %) to the new
# snapshot values, and fixes the row order / data for RenderToken vs
# ARBITRUM (rows 32 and 33 were swapped in the source feed).
#
# Price cells that look like plain numbers ("309.01", "1.007", ...) are
# written via a temporary "@" (text) number format so the COM layer keeps
# them as literal strings -- matching the source feed, which always
# stores these as text -- and then the cell style is reset back to
# "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.731.10"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.836.08"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4667"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3613"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07151"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9336"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.87%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07671"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "1.864.66"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.266"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.364"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008559"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "26.762.92"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.021"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.002"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.890"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08831"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.158"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.842"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.177"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7412"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.085"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.974"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01925"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05141"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5077"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1504"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.124"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4673"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06043"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  -0.31%  "
